$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on sheet "Hoja1" (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰" + "`n" + `
  "✅ Dólar paralelo: 68" + "`n" + `
  "" + "`n" + `
  "Binance" + "`n" + `
  "✅ 1000 Bs = 1.61 = 5967.5 pesos" + "`n" + `
  "✅ 5967.5 pesos = 1.63 = 970.43 Bs" + "`n" + `
  "" + "`n" + `
  "Promedio competencia" + "`n" + `
  "✅ Tasa pesos: 20" + "`n" + `
  "✅ Tasa Bs: 20" + "`n" + `
  "✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 620
$wsTasas.Range("O10").Value = 3699.85
$wsTasas.Range("N12").Value = 3665
$wsTasas.Range("O12").Value = 596
